# feat: add 2022-Q1 data
#
# - Inserts a new worksheet "2022-Q1" (fund-holdings snapshot) right before
#   the "总计" (totals) summary sheet.
# - Rebuilds "总计" with a new leading row for 2022-Q1 and the previously
#   existing rows shifted down by one.

$wb = $excel.ActiveWorkbook

function Set-HeaderCellStyle($cell) {
    # Matches the bold / centered / thin-bordered look used for the header
    # row and the leading index column ("s=2" in the original workbook).
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4160     # xlTop
    $cell.Borders.LineStyle = 1
}

function Set-TextValue($cell, $text) {
    # Forces a value to be stored as text (keeps leading zeros / avoids
    # numeric coercion of things like "001044" or "8.92").
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# ---------------------------------------------------------------------
# 1. Capture the existing "总计" data before we touch anything, then
#    delete the sheet so it can be rebuilt (with an extra row) at the
#    very end of the workbook, after the new "2022-Q1" sheet.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Delete()

# ---------------------------------------------------------------------
# 2. Add the new "2022-Q1" worksheet after the current last sheet
#    (i.e. "2021-Q4"), matching the fund-holdings layout used by the
#    other quarterly sheets.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$q1Sheet = $wb.Worksheets.Add($null, $lastSheet)
$q1Sheet.Name = "2022-Q1"

Set-HeaderCellStyle($q1Sheet.Range("B1"))
$q1Sheet.Range("B1").Value = "基金代码"
Set-HeaderCellStyle($q1Sheet.Range("C1"))
$q1Sheet.Range("C1").Value = "基金名称"
Set-HeaderCellStyle($q1Sheet.Range("D1"))
$q1Sheet.Range("D1").Value = "基金规模"
Set-HeaderCellStyle($q1Sheet.Range("E1"))
$q1Sheet.Range("E1").Value = "股票总仓位"
Set-HeaderCellStyle($q1Sheet.Range("F1"))
$q1Sheet.Range("F1").Value = "仓位占比"
Set-HeaderCellStyle($q1Sheet.Range("G1"))
$q1Sheet.Range("G1").Value = "持有市值(亿元)"
Set-HeaderCellStyle($q1Sheet.Range("H1"))
$q1Sheet.Range("H1").Value = "仓位排名"

Set-HeaderCellStyle($q1Sheet.Range("A2"))
$q1Sheet.Range("A2").Value = 0
Set-TextValue $q1Sheet.Range("B2") "001044"
Set-TextValue $q1Sheet.Range("C2") "嘉实新消费股票"
Set-TextValue $q1Sheet.Range("D2") "8.92"
Set-TextValue $q1Sheet.Range("E2") "80.25"
Set-TextValue $q1Sheet.Range("F2") "5.17"
Set-TextValue $q1Sheet.Range("G2") "0.4612"
$q1Sheet.Range("H2").Value = 8

# ---------------------------------------------------------------------
# 3. Rebuild "总计" after "2022-Q1", with the new 2022-Q1 row on top and
#    every previously existing row shifted down by one.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Add($null, $q1Sheet)
$totalSheet.Name = "总计"

Set-HeaderCellStyle($totalSheet.Range("B1"))
$totalSheet.Range("B1").Value = "日期"
Set-HeaderCellStyle($totalSheet.Range("C1"))
$totalSheet.Range("C1").Value = "持有数量(只)"
Set-HeaderCellStyle($totalSheet.Range("D1"))
$totalSheet.Range("D1").Value = "持有市值(亿元)"

$rows = @(
    @(0, "2022-Q1", 1, 0.46),
    @(1, "2021-Q4", 1, 0.49),
    @(2, "2021-Q2", 2, 0.45),
    @(3, "2021-Q1", 1, 0.05),
    @(4, "2020-Q4", 1, 0.08)
)

$r = 2
foreach ($row in $rows) {
    Set-HeaderCellStyle($totalSheet.Range("A$r"))
    $totalSheet.Range("A$r").Value = $row[0]
    $totalSheet.Range("B$r").Value = $row[1]
    $totalSheet.Range("C$r").Value = $row[2]
    $totalSheet.Range("D$r").Value = $row[3]
    $r = $r + 1
}
